# "1st changes of mifos to finflux"
#
# On the "Repayment Schedule" sheet, insert a new blank column before
# column N. This shifts the old N/O/P columns ("Late" header/data, the
# blank spacer column, and the "Outstanding" header/data) one column to
# the right into O/P/Q, and extends the used range from A1:P14 to
# A1:Q14.
#
# Also switch the active tab from "Transactions" to "Repayment Schedule"
# (with cell S5 selected there), which un-marks "Transactions" as the
# active tab (its own remembered selection, D10, is left untouched).

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts N,O,P -> O,P,Q).
$wsSchedule.Columns("N:N").Insert()

# Make "Repayment Schedule" the active sheet/tab, with S5 selected.
$wsSchedule.Activate()
$wsSchedule.Range("S5").Select()
